# OTreeAlg: confirm (SVS) instead of identify (ADS) for next states of extra states
#
# The underlying data sheet ("Moore_R100") used to keep a redundant column L
# with the whole row pre-joined as one semicolon-separated string
# (e.g. "1;456;912;5;2598;...;TeacherDFSM;;"). That column is retired here;
# the real, already-present columns A:J remain the single source of truth.
#
# On top of that, six new result rows are appended for the new "OTree"
# algorithm variant (with "ExtraStates: 1" and "ExtraStates:
# 1+confirm/identify" configurations), one row per Teacher
# (TeacherDFSM / TeacherRL / TeacherBB).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Moore_R100")

# --- Drop the redundant, pre-joined column L (rows 1-28 are the only ones
#     that still carry it) -----------------------------------------------
$ws.Range("L1:L28").ClearContents()

# --- Append the new OTree rows ------------------------------------------
$otreeRows = @(
    @(38, 1, 2798, 6337, 0, 20629, "../../data/tests/sequences/Moore_R100.fsm", "OTree", "ExtraStates: 1", "TeacherDFSM", ""),
    @(39, 1, 2798, 6337, 0, 20629, "../../data/tests/sequences/Moore_R100.fsm", "OTree", "ExtraStates: 1", "TeacherRL", ""),
    @(40, 1, 2798, 6337, 0, 6337, "../../data/tests/sequences/Moore_R100.fsm", "OTree", "ExtraStates: 1", "TeacherBB:SPY_method (3 extra states)", "BlackBoxDFSM"),
    @(41, 1, 2627, 5058, 0, 18344, "../../data/tests/sequences/Moore_R100.fsm", "OTree", "ExtraStates: 1+confirm/identify", "TeacherDFSM", ""),
    @(42, 1, 2627, 5058, 0, 18344, "../../data/tests/sequences/Moore_R100.fsm", "OTree", "ExtraStates: 1+confirm/identify", "TeacherRL", ""),
    @(43, 1, 2627, 5058, 0, 5058, "../../data/tests/sequences/Moore_R100.fsm", "OTree", "ExtraStates: 1+confirm/identify", "TeacherBB:SPY_method (3 extra states)", "BlackBoxDFSM")
)

foreach ($r in $otreeRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
    $ws.Cells.Item($row, 5).Value = $r[5]
    $ws.Cells.Item($row, 6).Value = $r[6]
    $ws.Cells.Item($row, 7).Value = $r[7]
    $ws.Cells.Item($row, 8).Value = $r[8]
    $ws.Cells.Item($row, 9).Value = $r[9]
    if ($r[10] -ne "") {
        $ws.Cells.Item($row, 10).Value = $r[10]
    }
}

# --- Columns C and D now share the same ("best fit") width ---------------
$ws.Columns.Item(3).ColumnWidth = 5.14
$ws.Columns.Item(4).ColumnWidth = 5.14

# --- Reflect the new scroll position / selection used after the edit -----
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 30
$win.ScrollColumn = 1
$ws.Range("B44:H48").Select()
